$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new blank rows above the current row 2 (old data shifts from
# rows 2-10 down to rows 6-14).
$ws.Range("A2:A5").EntireRow.Insert()

# --- New row 2: "ItemMasterDetails, Lookup multple value, Item Maste Add Edit single forms"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "ItemMasterDetails, Lookup multple value, Item Maste Add Edit single forms"
$ws.Range("E2").Value = "NO"
$ws.Range("F2").Value = "In progress"

# --- New row 3: "Party master"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Party master"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "In Progress"

# --- New row 4: "Multple selection in Transation forms"
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Multple selection in Transation forms"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = " In Progress"

# Row 5 stays completely empty (gap row).

# --- Add a new trailing row 15 with just the next sequence number.
$ws.Range("A15").Value = 10

# --- Row 14 (previously row 10) gains a "pending in reports" note in G.
$ws.Range("G14").Value = "pending in reports"

# Refresh the view/selection + used-range dimension to match the edit.
$ws.Range("E5").Select() | Out-Null
